# Correct model name: replace "B9 Tribeca (2006 - 2007)\nTribeca (2008 - onward)"
# with just "Tribeca" in column D for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(18, 21, 29, 31, 36, 44, 46, 52)

foreach ($r in $rows) {
    $cell = $ws.Range("D$r")
    $current = $cell.Value()
    if ($current -like "B9 Tribeca*") {
        $cell.Value = "Tribeca"
    }
}
